$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 501.5
$ws.Range("I2").Value = 448.54544
$ws.Range("J2").Value = 584.7143
$ws.Range("K2").Value = 448.54544
$ws.Range("L2").Value = 584.7143
$ws.Range("M2").Value = -335.54544
$ws.Range("N2").Value = -810.7143

$ws.Range("H88").Value = 1587
$ws.Range("I88").Value = 3817.3333
$ws.Range("J88").Value = 750.625
$ws.Range("K88").Value = 3817.3333
$ws.Range("L88").Value = 750.625
$ws.Range("M88").Value = -3411.3333
$ws.Range("N88").Value = -1562.625

$ws.Range("H91").Value = 1587
$ws.Range("I91").Value = 3817.3333
$ws.Range("J91").Value = 750.625
$ws.Range("K91").Value = 3817.3333
$ws.Range("L91").Value = 750.625
$ws.Range("M91").Value = -2413.3333
$ws.Range("N91").Value = -3558.625

$ws.Range("H100").Value = 1618.5264
$ws.Range("J100").Value = 3999.8
$ws.Range("L100").Value = 3999.8
$ws.Range("N100").Value = -5081.8

$ws.Range("H137").Value = 76925830
$ws.Range("I137").Value = 45457440
$ws.Range("K137").Value = 136372320
$ws.Range("M137").Value = -136369770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10646924
$ws.Range("I32").Value = 13892652
$ws.Range("J32").Value = 24539.773
$ws.Range("K32").Value = 13892652
$ws.Range("L32").Value = 24539.773
$ws.Range("M32").Value = -13892365
$ws.Range("N32").Value = -25113.773

$ws.Range("H74").Value = 71510960
$ws.Range("I74").Value = 77011540
$ws.Range("K74").Value = 77011540
$ws.Range("M74").Value = -77010666

$ws.Range("H77").Value = 71510960
$ws.Range("I77").Value = 77011540
$ws.Range("K77").Value = 385057700
$ws.Range("M77").Value = -385053332

$ws.Range("H80").Value = 77331.75
$ws.Range("J80").Value = 79554
$ws.Range("L80").Value = 79554
$ws.Range("N80").Value = -81550

$ws.Range("H83").Value = 77331.75
$ws.Range("J83").Value = 79554
$ws.Range("L83").Value = 238662
$ws.Range("N83").Value = -248646

$ws.Range("H97").Value = 1018.6667
$ws.Range("I97").Value = 427.25
$ws.Range("K97").Value = 427.25
$ws.Range("M97").Value = 68.75

$ws.Range("H110").Value = 15708.179
$ws.Range("I110").Value = 17627.652
$ws.Range("J110").Value = 6878.6
$ws.Range("K110").Value = 17627.652
$ws.Range("L110").Value = 6878.6
$ws.Range("M110").Value = -15582.652
$ws.Range("N110").Value = -10968.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1250
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -1726

$ws.Range("H86").Value = 20839.723
$ws.Range("I86").Value = 13580.5
$ws.Range("J86").Value = 29913.75
$ws.Range("K86").Value = 13580.5
$ws.Range("L86").Value = 29913.75
$ws.Range("M86").Value = -12457.5
$ws.Range("N86").Value = -32159.75

$ws.Range("H89").Value = 20839.723
$ws.Range("I89").Value = 13580.5
$ws.Range("J89").Value = 29913.75
$ws.Range("K89").Value = 67902.5
$ws.Range("L89").Value = 149568.75
$ws.Range("M89").Value = -62286.5
$ws.Range("N89").Value = -160800.75

$ws.Range("H107").Value = 7583.591
$ws.Range("I107").Value = 6066.6665
$ws.Range("J107").Value = 8152.4375
$ws.Range("K107").Value = 6066.6665
$ws.Range("L107").Value = 8152.4375
$ws.Range("M107").Value = -4146.6665
$ws.Range("N107").Value = -11992.4375

$ws.Range("H134").Value = 4323.5713
$ws.Range("I134").Value = 2862.6785
$ws.Range("J134").Value = 10167.143
$ws.Range("K134").Value = 8588.0355
$ws.Range("L134").Value = 30501.429
$ws.Range("M134").Value = -6053.0355
$ws.Range("N134").Value = -35571.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1826.6428
$ws.Range("J16").Value = 990.5
$ws.Range("L16").Value = 990.5
$ws.Range("N16").Value = -1564.5

$ws.Range("H31").Value = 27031462
$ws.Range("I31").Value = 3548.5518
$ws.Range("K31").Value = 3548.5518
$ws.Range("M31").Value = -3253.5518

$ws.Range("H34").Value = 27031462
$ws.Range("I34").Value = 3548.5518
$ws.Range("K34").Value = 3548.5518
$ws.Range("M34").Value = -3346.5518

$ws.Range("H58").Value = 5333.5557
$ws.Range("I58").Value = 5875.875
$ws.Range("J58").Value = 995
$ws.Range("K58").Value = 5875.875
$ws.Range("L58").Value = 995
$ws.Range("M58").Value = -5672.875
$ws.Range("N58").Value = -1401

$ws.Range("H107").Value = 1643.9546
$ws.Range("I107").Value = 919.5
$ws.Range("J107").Value = 2247.6667
$ws.Range("K107").Value = 919.5
$ws.Range("L107").Value = 2247.6667
$ws.Range("M107").Value = 1000.5
$ws.Range("N107").Value = -6087.6667

$ws.Range("H113").Value = 1826.6428
$ws.Range("J113").Value = 990.5
$ws.Range("L113").Value = 990.5
$ws.Range("N113").Value = -5330.5

$ws.Range("H122").Value = 2258.6365
$ws.Range("I122").Value = 2264.75
$ws.Range("J122").Value = 2255.1428
$ws.Range("K122").Value = 6794.25
$ws.Range("L122").Value = 6765.428400000001
$ws.Range("M122").Value = -4344.25
$ws.Range("N122").Value = -11665.4284

$ws.Range("H132").Value = 85857.03999999999
$ws.Range("I132").Value = 101867.05
$ws.Range("K132").Value = 305601.15
$ws.Range("M132").Value = -303071.15

$ws.Range("H136").Value = 5333.5557
$ws.Range("I136").Value = 5875.875
$ws.Range("J136").Value = 995
$ws.Range("K136").Value = 17627.625
$ws.Range("L136").Value = 2985
$ws.Range("M136").Value = -15077.625
$ws.Range("N136").Value = -8085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1372.5
$ws.Range("J25").Value = 3011.25
$ws.Range("L25").Value = 9033.75
$ws.Range("N25").Value = -9371.75

$ws.Range("H30").Value = 1372.5
$ws.Range("J30").Value = 3011.25
$ws.Range("L30").Value = 9033.75
$ws.Range("N30").Value = -9237.75

$ws.Range("H34").Value = 2331.5
$ws.Range("J34").Value = 3247.5
$ws.Range("L34").Value = 9742.5
$ws.Range("N34").Value = -9910.5

$ws.Range("H80").Value = 4333.3335
$ws.Range("J80").Value = 4333.3335
$ws.Range("L80").Value = 13000.0005
$ws.Range("N80").Value = -14872.0005

$ws.Range("H83").Value = 4333.3335
$ws.Range("J83").Value = 4333.3335
$ws.Range("L83").Value = 39000.0015
$ws.Range("N83").Value = -48360.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4335.2856
$ws.Range("I80").Value = 4058.25
$ws.Range("K80").Value = 4058.25
$ws.Range("M80").Value = -3060.25

$ws.Range("H83").Value = 4335.2856
$ws.Range("I83").Value = 4058.25
$ws.Range("K83").Value = 20291.25
$ws.Range("M83").Value = -15299.25

$ws.Range("H97").Value = 1260
$ws.Range("I97").Value = 724.0714
$ws.Range("J97").Value = 1837.1538
$ws.Range("K97").Value = 724.0714
$ws.Range("L97").Value = 1837.1538
$ws.Range("M97").Value = -228.0714
$ws.Range("N97").Value = -2829.1538

$ws.Range("H113").Value = 4172.643
$ws.Range("I113").Value = 3554.6667
$ws.Range("J113").Value = 5285
$ws.Range("K113").Value = 3554.6667
$ws.Range("L113").Value = 5285
$ws.Range("M113").Value = -1384.6667
$ws.Range("N113").Value = -9625

$ws.Range("H127").Value = 99999
$ws.Range("J127").Value = 99999
$ws.Range("L127").Value = 99999
$ws.Range("N127").Value = -109919

$ws.Range("H132").Value = 1918.8654
$ws.Range("I132").Value = 1699
$ws.Range("K132").Value = 5097
$ws.Range("M132").Value = -2567

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1287.9546
$ws.Range("I46").Value = 899.8889
$ws.Range("K46").Value = 899.8889
$ws.Range("M46").Value = -711.8889

$ws.Range("H82").Value = 3427.389
$ws.Range("I82").Value = 1299
$ws.Range("J82").Value = 5555.778
$ws.Range("K82").Value = 1299
$ws.Range("L82").Value = 5555.778
$ws.Range("M82").Value = -938
$ws.Range("N82").Value = -6277.778

$ws.Range("H85").Value = 3427.389
$ws.Range("I85").Value = 1299
$ws.Range("J85").Value = 5555.778
$ws.Range("K85").Value = 1299
$ws.Range("L85").Value = 5555.778
$ws.Range("M85").Value = -51
$ws.Range("N85").Value = -8051.778

$ws.Range("H93").Value = 2103.75
$ws.Range("I93").Value = 1155.25
$ws.Range("J93").Value = 4949.25
$ws.Range("K93").Value = 1155.25
$ws.Range("L93").Value = 4949.25
$ws.Range("M93").Value = 92.75
$ws.Range("N93").Value = -7445.25

$ws.Range("H100").Value = 3197.0667
$ws.Range("I100").Value = 2271.8572
$ws.Range("J100").Value = 4006.625
$ws.Range("K100").Value = 2271.8572
$ws.Range("L100").Value = 4006.625
$ws.Range("M100").Value = -1730.8572
$ws.Range("N100").Value = -5088.625

$ws.Range("H136").Value = 3164.4
$ws.Range("I136").Value = 3214.8975
$ws.Range("J136").Value = 1195
$ws.Range("K136").Value = 9644.692500000001
$ws.Range("L136").Value = 3585
$ws.Range("M136").Value = -7094.692500000001
$ws.Range("N136").Value = -8685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 10002503
$ws.Range("J10").Value = 10002503
$ws.Range("L10").Value = 10002503
$ws.Range("N10").Value = -10002841

$ws.Range("H132").Value = 5974.1387
$ws.Range("I132").Value = 5579.8286
$ws.Range("J132").Value = 19775
$ws.Range("K132").Value = 16739.4858
$ws.Range("L132").Value = 59325
$ws.Range("M132").Value = -14209.4858
$ws.Range("N132").Value = -64385

Write-Host "Applied profit-sheet updates"